$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "1f4c7be5-410a-42c2-a297-a89e99fde061"
$ws.Range("A9").Value = "c96e1fd1-f92f-422c-b339-dd4223320a26"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "29"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "39"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "54"
